$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.057.85'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').Value = '2.558.16'
$ws.Range('E3').Value = '  +0.86%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '580.69'
$ws.Range('E5').Value = '  +2.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.06'
$ws.Range('E6').Value = '  -0.68%  '
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('E8').Value = '  +0.73%  '
$ws.Range('E9').Value = '  +0.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.60'
$ws.Range('E10').Value = '  -0.63%  '
$ws.Range('E11').Value = '  +0.09%  '
$ws.Range('E12').Value = '  +0.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.49'
$ws.Range('E13').Value = '  -0.95%  '
$ws.Range('D14').Value = '3.018.34'
$ws.Range('E14').Value = '  +0.61%  '
$ws.Range('D15').Value = '62.984.16'
$ws.Range('E15').Value = '  +0.01%  '
$ws.Range('E16').Value = '  +2.03%  '
$ws.Range('D17').Value = '2.565.87'
$ws.Range('E17').Value = '  +0.67%  '
$ws.Range('E18').Value = '  -1.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '340.78'
$ws.Range('E19').Value = '  +1.27%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.36'
$ws.Range('E20').Value = '  +2.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.81'
$ws.Range('E21').Value = '  +1.38%  '
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '65.88'
$ws.Range('E23').Value = '  +0.24%  '
$ws.Range('D24').Value = '2.678.20'
$ws.Range('E24').Value = '  +0.17%  '
$ws.Range('B25').Value = 'Fetch.AI'
$ws.Range('C25').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.63'
$ws.Range('E25').Value = '  +3.61%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.170'
$ws.Range('E26').Value = '  +1.30%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('B28').Value = 'SuiNetwork'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.49'
$ws.Range('E28').Value = '  -1.25%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.41'
$ws.Range('E29').Value = '  +1.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.88'
$ws.Range('E30').Value = '  +9.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.96'
$ws.Range('E31').Value = '  +5.78%  '
$ws.Range('D32').Value = '0.0₃0819'
$ws.Range('E32').Value = '  +1.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '177.41'
$ws.Range('E33').Value = '  -0.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.59'
$ws.Range('E34').Value = '  +1.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '425.61'
$ws.Range('E35').Value = '  +2.68%  '
$ws.Range('E36').Value = '  +1.26%  '
$ws.Range('E37').Value = '  +2.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.43'
$ws.Range('E38').Value = '  +1.31%  '
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.69'
$ws.Range('E42').Value = '  +0.81%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '151.37'
$ws.Range('E43').Value = '  -0.48%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.80'
$ws.Range('E44').Value = '  +1.76%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '20.88'
$ws.Range('E45').Value = '  +1.71%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0547'
$ws.Range('E46').Value = '  +5.49%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.605'
$ws.Range('E47').Value = '  -0.12%  '
$ws.Range('E48').Value = '  +1.31%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0241'
$ws.Range('E49').Value = '  +1.92%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.35'
$ws.Range('E50').Value = '  -0.08%  '
$ws.Range('E51').Value = '  -2.88%  '
